# Remove `pax` from databases
# Rename use-type property headers that used a "...pax" suffix to the
# new "...p" style naming convention, on both sheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- INTERNAL_LOADS sheet ---
$ws2 = $wb.Worksheets.Item("INTERNAL_LOADS")
$ws2.Range("B1").Value = "Occ_m2p"
$ws2.Range("C1").Value = "Qs_Wp"
$ws2.Range("D1").Value = "X_ghp"
$ws2.Range("I1").Value = "Vww_ldp"
$ws2.Range("J1").Value = "Vw_ldp"

# --- INDOOR_COMFORT sheet ---
$ws1 = $wb.Worksheets.Item("INDOOR_COMFORT")
$ws1.Range("F1").Value = "Ve_lsp"

# Make INDOOR_COMFORT the active/selected sheet, matching the commit.
$ws1.Activate()
$ws1.Range("F2").Select()
